$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The eccentricity equation ($VALUE@eccentricity@EQUATIONS, column I) is being
# changed for every configuration row: it no longer rounds the computed value,
# allowing the two cycloidal discs to be offset asymmetrically so the gearbox
# meshes correctly.
$newEccentricityFormula = "'=(""pin_pitch_dia"" / 2) / ""N_pins"" * 0.7"

$ws.Range("I2").Value = $newEccentricityFormula
$ws.Range("I3").Value = $newEccentricityFormula
$ws.Range("I4").Value = $newEccentricityFormula
$ws.Range("I5").Value = $newEccentricityFormula
